$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Depart On" (E2) / "Return On" (F2) / E3 previously held raw date serials;
# replace them with pre-formatted date-text strings (new shared strings).
$ws.Range("E2").Value = "Thu, 14 Feb, 2019"
$ws.Range("F2").Value = "Sat, 15 Feb, 2019"
$ws.Range("E3").Value = "Thu, 16 Sun, 2019"

# Column E needs to widen to fit the longer text now stored in it.
$ws.Columns.Item(5).ColumnWidth = 15.3

# Move the active selection from J3 to E3.
$ws.Range("E3").Select()
